$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 935.2  # H28 was 936.8570999999999
$ws.Cells.Item(28, 9).Value = 816.8889  # I28 was 936.8570999999999
$ws.Cells.Item(28, 10).Value = 2000  # J28 was 0
$ws.Cells.Item(28, 11).Value = 816.8889  # K28 was 936.8570999999999
$ws.Cells.Item(28, 12).Value = 2000  # L28 was 0
$ws.Cells.Item(28, 13).Value = -331.8889  # M28 was -451.8570999999999
$ws.Cells.Item(28, 14).Value = -2970  # N28 was None
$ws.Cells.Item(96, 8).Value = 1170.4  # H96 was 1109.6875
$ws.Cells.Item(96, 9).Value = 605.8  # I96 was 568.8182
$ws.Cells.Item(96, 11).Value = 1817.4  # K96 was 1706.4546
$ws.Cells.Item(96, 13).Value = -444.3999999999999  # M96 was -333.4546
$ws.Cells.Item(100, 8).Value = 3900.0833  # H100 was 4027.3635
$ws.Cells.Item(100, 10).Value = 2499.5  # J100 was 2499
$ws.Cells.Item(100, 12).Value = 2499.5  # L100 was 2499
$ws.Cells.Item(100, 14).Value = -3581.5  # N100 was -3581
$ws.Cells.Item(107, 8).Value = 679.5714  # H107 was 694.26666
$ws.Cells.Item(107, 10).Value = 1117  # J107 was 1062.75
$ws.Cells.Item(107, 12).Value = 1117  # L107 was 1062.75
$ws.Cells.Item(107, 14).Value = -4957  # N107 was -4902.75
$ws.Cells.Item(113, 8).Value = 8860.75  # H113 was 10814.667
$ws.Cells.Item(113, 9).Value = 2999.5  # I113 was 3000
$ws.Cells.Item(113, 11).Value = 2999.5  # K113 was 3000
$ws.Cells.Item(113, 13).Value = 254.5  # M113 was 254
$ws.Cells.Item(132, 8).Value = 3336  # H132 was 3659.6
$ws.Cells.Item(132, 9).Value = 2260.923  # I132 was 2635.6365
$ws.Cells.Item(132, 10).Value = 4888.8887  # J132 was 4911.1113
$ws.Cells.Item(132, 11).Value = 6782.768999999999  # K132 was 7906.9095
$ws.Cells.Item(132, 12).Value = 14666.6661  # L132 was 14733.3339
$ws.Cells.Item(132, 13).Value = -4252.768999999999  # M132 was -5376.9095
$ws.Cells.Item(132, 14).Value = -19726.6661  # N132 was -19793.3339
$ws.Cells.Item(137, 8).Value = 2807.92  # H137 was 2638.0952
$ws.Cells.Item(137, 9).Value = 1383.3077  # I137 was 1432.0834
$ws.Cells.Item(137, 10).Value = 4351.25  # J137 was 4246.1113
$ws.Cells.Item(137, 11).Value = 4149.9231  # K137 was 4296.2502
$ws.Cells.Item(137, 12).Value = 13053.75  # L137 was 12738.3339
$ws.Cells.Item(137, 13).Value = -1599.9231  # M137 was -1746.2502
$ws.Cells.Item(137, 14).Value = -18153.75  # N137 was -17838.3339
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1066.25  # H2 was 1072.5
$ws.Cells.Item(2, 9).Value = 1066.25  # I2 was 1072.5
$ws.Cells.Item(2, 11).Value = 1066.25  # K2 was 1072.5
$ws.Cells.Item(2, 13).Value = -953.25  # M2 was -959.5
$ws.Cells.Item(18, 8).Value = 1999  # H18 was 0
$ws.Cells.Item(18, 10).Value = 1999  # J18 was 0
$ws.Cells.Item(18, 12).Value = 1999  # L18 was 0
$ws.Cells.Item(18, 14).Value = -2643  # N18 was None
$ws.Cells.Item(45, 8).Value = 2385.625  # H45 was 2999.5715
$ws.Cells.Item(45, 9).Value = 2027.75  # I45 was 2333.3333
$ws.Cells.Item(45, 10).Value = 2743.5  # J45 was 3499.25
$ws.Cells.Item(45, 11).Value = 2027.75  # K45 was 2333.3333
$ws.Cells.Item(45, 12).Value = 2743.5  # L45 was 3499.25
$ws.Cells.Item(45, 13).Value = -1650.75  # M45 was -1956.3333
$ws.Cells.Item(45, 14).Value = -3497.5  # N45 was -4253.25
$ws.Cells.Item(61, 8).Value = 5410.6665  # H61 was 4969.6
$ws.Cells.Item(61, 9).Value = 1399.3334  # I61 was 1299.5
$ws.Cells.Item(61, 11).Value = 1399.3334  # K61 was 1299.5
$ws.Cells.Item(61, 13).Value = -1187.3334  # M61 was -1087.5
$ws.Cells.Item(110, 8).Value = 2906.7273  # H110 was 3186.111
$ws.Cells.Item(110, 9).Value = 1570.7142  # I110 was 1782.6666
$ws.Cells.Item(110, 10).Value = 5244.75  # J110 was 5993
$ws.Cells.Item(110, 11).Value = 1570.7142  # K110 was 1782.6666
$ws.Cells.Item(110, 12).Value = 5244.75  # L110 was 5993
$ws.Cells.Item(110, 13).Value = 474.2858000000001  # M110 was 262.3334
$ws.Cells.Item(110, 14).Value = -9334.75  # N110 was -10083
$ws.Cells.Item(116, 8).Value = 1066.25  # H116 was 1072.5
$ws.Cells.Item(116, 9).Value = 1066.25  # I116 was 1072.5
$ws.Cells.Item(116, 11).Value = 1066.25  # K116 was 1072.5
$ws.Cells.Item(116, 13).Value = 1227.75  # M116 was 1221.5
$ws.Cells.Item(124, 8).Value = 45000  # H124 was 44933.332
$ws.Cells.Item(124, 10).Value = 45000  # J124 was 44933.332
$ws.Cells.Item(124, 12).Value = 45000  # L124 was 44933.332
$ws.Cells.Item(124, 14).Value = -54820  # N124 was -54753.332
$ws.Cells.Item(125, 8).Value = 0  # H125 was 60000
$ws.Cells.Item(125, 9).Value = 0  # I125 was 60000
$ws.Cells.Item(125, 11).Value = 0  # K125 was 60000
$ws.Cells.Item(125, 13).ClearContents()  # M125 was -55080
$ws.Cells.Item(136, 8).Value = 5410.6665  # H136 was 4969.6
$ws.Cells.Item(136, 9).Value = 1399.3334  # I136 was 1299.5
$ws.Cells.Item(136, 11).Value = 4198.0002  # K136 was 3898.5
$ws.Cells.Item(136, 13).Value = -1648.0002  # M136 was -1348.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1066.25  # H3 was 1072.5
$ws.Cells.Item(3, 9).Value = 1066.25  # I3 was 1072.5
$ws.Cells.Item(3, 11).Value = 1066.25  # K3 was 1072.5
$ws.Cells.Item(3, 13).Value = -952.25  # M3 was -958.5
$ws.Cells.Item(20, 8).Value = 5245.5  # H20 was 5926.2
$ws.Cells.Item(20, 9).Value = 2920.8572  # I20 was 3648.8
$ws.Cells.Item(20, 10).Value = 10669.667  # J20 was 8203.6
$ws.Cells.Item(20, 11).Value = 2920.8572  # K20 was 3648.8
$ws.Cells.Item(20, 12).Value = 10669.667  # L20 was 8203.6
$ws.Cells.Item(20, 13).Value = -2673.8572  # M20 was -3401.8
$ws.Cells.Item(20, 14).Value = -11163.667  # N20 was -8697.6
$ws.Cells.Item(82, 8).Value = 11510.667  # H82 was 13859.667
$ws.Cells.Item(82, 9).Value = 11510.667  # I82 was 13859.667
$ws.Cells.Item(82, 11).Value = 11510.667  # K82 was 13859.667
$ws.Cells.Item(82, 13).Value = -11127.667  # M82 was -13476.667
$ws.Cells.Item(85, 8).Value = 11510.667  # H85 was 13859.667
$ws.Cells.Item(85, 9).Value = 11510.667  # I85 was 13859.667
$ws.Cells.Item(85, 11).Value = 11510.667  # K85 was 13859.667
$ws.Cells.Item(85, 13).Value = -10184.667  # M85 was -12533.667
$ws.Cells.Item(86, 8).Value = 2665.6667  # H86 was 3250
$ws.Cells.Item(86, 10).Value = 3248.5  # J86 was 5000
$ws.Cells.Item(86, 12).Value = 3248.5  # L86 was 5000
$ws.Cells.Item(86, 14).Value = -5494.5  # N86 was -7246
$ws.Cells.Item(89, 8).Value = 2665.6667  # H89 was 3250
$ws.Cells.Item(89, 10).Value = 3248.5  # J89 was 5000
$ws.Cells.Item(89, 12).Value = 16242.5  # L89 was 25000
$ws.Cells.Item(89, 14).Value = -27474.5  # N89 was -36232
$ws.Cells.Item(105, 8).Value = 4638.643  # H105 was 4820.7144
$ws.Cells.Item(105, 9).Value = 2490.3333  # I105 was 2915.1667
$ws.Cells.Item(105, 11).Value = 2490.3333  # K105 was 2915.1667
$ws.Cells.Item(105, 13).Value = -743.3332999999998  # M105 was -1168.1667
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 9516.666999999999  # H4 was 11219
$ws.Cells.Item(4, 9).Value = 6800  # I4 was 8497.5
$ws.Cells.Item(4, 10).Value = 11457.143  # J4 was 13033.333
$ws.Cells.Item(4, 11).Value = 6800  # K4 was 8497.5
$ws.Cells.Item(4, 12).Value = 11457.143  # L4 was 13033.333
$ws.Cells.Item(4, 13).Value = -6688  # M4 was -8385.5
$ws.Cells.Item(4, 14).Value = -11681.143  # N4 was -13257.333
$ws.Cells.Item(31, 8).Value = 3087.6924  # H31 was 2836.8
$ws.Cells.Item(31, 9).Value = 2970.1667  # I31 was 2718.1428
$ws.Cells.Item(31, 11).Value = 2970.1667  # K31 was 2718.1428
$ws.Cells.Item(31, 13).Value = -2675.1667  # M31 was -2423.1428
$ws.Cells.Item(32, 8).Value = 3479  # H32 was 5750
$ws.Cells.Item(32, 9).Value = 3848.75  # I32 was 5750
$ws.Cells.Item(32, 10).Value = 2000  # J32 was 0
$ws.Cells.Item(32, 11).Value = 3848.75  # K32 was 5750
$ws.Cells.Item(32, 12).Value = 2000  # L32 was 0
$ws.Cells.Item(32, 13).Value = -3532.75  # M32 was -5434
$ws.Cells.Item(32, 14).Value = -2632  # N32 was None
$ws.Cells.Item(34, 8).Value = 3087.6924  # H34 was 2836.8
$ws.Cells.Item(34, 9).Value = 2970.1667  # I34 was 2718.1428
$ws.Cells.Item(34, 11).Value = 2970.1667  # K34 was 2718.1428
$ws.Cells.Item(34, 13).Value = -2768.1667  # M34 was -2516.1428
$ws.Cells.Item(50, 8).Value = 65000  # H50 was 47499
$ws.Cells.Item(50, 10).Value = 65000  # J50 was 47499
$ws.Cells.Item(50, 12).Value = 65000  # L50 was 47499
$ws.Cells.Item(50, 14).Value = -66250  # N50 was -48749
$ws.Cells.Item(92, 8).Value = 38000  # H92 was 32183.666
$ws.Cells.Item(92, 10).Value = 38000  # J92 was 32183.666
$ws.Cells.Item(92, 12).Value = 38000  # L92 was 32183.666
$ws.Cells.Item(92, 14).Value = -42992  # N92 was -37175.666
$ws.Cells.Item(94, 8).Value = 4265.4443  # H94 was 4550.5
$ws.Cells.Item(94, 9).Value = 4162.5  # I94 was 4598
$ws.Cells.Item(94, 11).Value = 4162.5  # K94 was 4598
$ws.Cells.Item(94, 13).Value = -3711.5  # M94 was -4147
$ws.Cells.Item(134, 8).Value = 3230.5186  # H134 was 3074.8965
$ws.Cells.Item(134, 9).Value = 3010.3044  # I134 was 2847.4
$ws.Cells.Item(134, 11).Value = 9030.913199999999  # K134 was 8542.200000000001
$ws.Cells.Item(134, 13).Value = -6495.913199999999  # M134 was -6007.200000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 44.333332  # H2 was 44.5
$ws.Cells.Item(2, 9).Value = 44.5  # I2 was 45
$ws.Cells.Item(2, 11).Value = 267  # K2 was 270
$ws.Cells.Item(2, 13).Value = -154  # M2 was -157
$ws.Cells.Item(12, 8).Value = 785.8333  # H12 was 913.3
$ws.Cells.Item(12, 9).Value = 181.14285  # I12 was 188.57143
$ws.Cells.Item(12, 10).Value = 1170.6364  # J12 was 1303.5385
$ws.Cells.Item(12, 11).Value = 543.4285500000001  # K12 was 565.71429
$ws.Cells.Item(12, 12).Value = 3511.9092  # L12 was 3910.6155
$ws.Cells.Item(12, 13).Value = -370.4285500000001  # M12 was -392.71429
$ws.Cells.Item(12, 14).Value = -3857.9092  # N12 was -4256.6155
$ws.Cells.Item(46, 8).Value = 151.5  # H46 was 184.75
$ws.Cells.Item(46, 9).Value = 161.8  # I46 was 213
$ws.Cells.Item(46, 11).Value = 485.4  # K46 was 639
$ws.Cells.Item(46, 13).Value = -394.4  # M46 was -548
$ws.Cells.Item(59, 8).Value = 1200  # H59 was 0
$ws.Cells.Item(59, 10).Value = 1200  # J59 was 0
$ws.Cells.Item(59, 12).Value = 3600  # L59 was 0
$ws.Cells.Item(59, 14).Value = -4680  # N59 was None
$ws.Cells.Item(69, 8).Value = 2471.111  # H69 was 488
$ws.Cells.Item(69, 9).Value = 2686.25  # I69 was 422.5
$ws.Cells.Item(69, 11).Value = 8058.75  # K69 was 1267.5
$ws.Cells.Item(69, 13).Value = -7247.75  # M69 was -456.5
$ws.Cells.Item(72, 8).Value = 2471.111  # H72 was 488
$ws.Cells.Item(72, 9).Value = 2686.25  # I72 was 422.5
$ws.Cells.Item(72, 11).Value = 24176.25  # K72 was 3802.5
$ws.Cells.Item(72, 13).Value = -20120.25  # M72 was 253.5
$ws.Cells.Item(111, 8).Value = 0  # H111 was 1000
$ws.Cells.Item(111, 9).Value = 0  # I111 was 1000
$ws.Cells.Item(111, 11).Value = 0  # K111 was 3000
$ws.Cells.Item(111, 13).ClearContents()  # M111 was 67
$ws.Cells.Item(114, 8).Value = 353.75  # H114 was 310.7
$ws.Cells.Item(114, 9).Value = 353.75  # I114 was 339.66666
$ws.Cells.Item(114, 10).Value = 0  # J114 was 50
$ws.Cells.Item(114, 11).Value = 1061.25  # K114 was 1018.99998
$ws.Cells.Item(114, 12).Value = 0  # L114 was 150
$ws.Cells.Item(114, 13).Value = 2192.75  # M114 was 2235.00002
$ws.Cells.Item(114, 14).ClearContents()  # N114 was -6658
$ws.Cells.Item(122, 8).Value = 1352.3334  # H122 was 1352.4667
$ws.Cells.Item(122, 10).Value = 1256.3636  # J122 was 1256.5454
$ws.Cells.Item(122, 12).Value = 11307.2724  # L122 was 11308.9086
$ws.Cells.Item(122, 14).Value = -16207.2724  # N122 was -16208.9086
$ws.Cells.Item(131, 8).Value = 1010.625  # H131 was 1124.2222
$ws.Cells.Item(131, 10).Value = 1997  # J131 was 2015
$ws.Cells.Item(131, 12).Value = 5991  # L131 was 6045
$ws.Cells.Item(131, 14).Value = -16071  # N131 was -16125
$ws.Cells.Item(139, 8).Value = 2950.2727  # H139 was 3542.077
$ws.Cells.Item(139, 9).Value = 1779.8572  # I139 was 1843.1666
$ws.Cells.Item(139, 10).Value = 4998.5  # J139 was 4998.2856
$ws.Cells.Item(139, 11).Value = 5339.571599999999  # K139 was 5529.4998
$ws.Cells.Item(139, 12).Value = 14995.5  # L139 was 14994.8568
$ws.Cells.Item(139, 13).Value = -199.5715999999993  # M139 was -389.4997999999996
$ws.Cells.Item(139, 14).Value = -25275.5  # N139 was -25274.8568
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 759.8  # H107 was 778.4
$ws.Cells.Item(107, 9).Value = 759.8  # I107 was 778.4
$ws.Cells.Item(107, 11).Value = 759.8  # K107 was 778.4
$ws.Cells.Item(107, 13).Value = 1160.2  # M107 was 1141.6
$ws.Cells.Item(113, 8).Value = 2357.1428  # H113 was 2533.2222
$ws.Cells.Item(113, 10).Value = 4500  # J113 was 3824.75
$ws.Cells.Item(113, 12).Value = 4500  # L113 was 3824.75
$ws.Cells.Item(113, 14).Value = -8840  # N113 was -8164.75
$ws.Cells.Item(132, 8).Value = 4658.7144  # H132 was 5142.4
$ws.Cells.Item(132, 9).Value = 4122.2  # I132 was 4178
$ws.Cells.Item(132, 10).Value = 6000  # J132 was 9000
$ws.Cells.Item(132, 11).Value = 12366.6  # K132 was 12534
$ws.Cells.Item(132, 12).Value = 18000  # L132 was 27000
$ws.Cells.Item(132, 13).Value = -9836.599999999999  # M132 was -10004
$ws.Cells.Item(132, 14).Value = -23060  # N132 was -32060
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 15000  # H2 was 3680
$ws.Cells.Item(2, 9).Value = 0  # I2 was 850
$ws.Cells.Item(2, 11).Value = 0  # K2 was 850
$ws.Cells.Item(2, 13).ClearContents()  # M2 was -738
$ws.Cells.Item(7, 8).Value = 1000  # H7 was 500
$ws.Cells.Item(7, 9).Value = 1000  # I7 was 500
$ws.Cells.Item(7, 11).Value = 1000  # K7 was 500
$ws.Cells.Item(7, 13).Value = -888  # M7 was -388
$ws.Cells.Item(11, 8).Value = 19999  # H11 was 0
$ws.Cells.Item(11, 10).Value = 19999  # J11 was 0
$ws.Cells.Item(11, 12).Value = 19999  # L11 was 0
$ws.Cells.Item(11, 14).Value = -20279  # N11 was None
$ws.Cells.Item(40, 8).Value = 13304.723  # H40 was 16123.75
$ws.Cells.Item(40, 9).Value = 13729.692  # I40 was 17185.75
$ws.Cells.Item(40, 10).Value = 12199.8  # J40 was 13999.75
$ws.Cells.Item(40, 11).Value = 13729.692  # K40 was 17185.75
$ws.Cells.Item(40, 12).Value = 12199.8  # L40 was 13999.75
$ws.Cells.Item(40, 13).Value = -13593.692  # M40 was -17049.75
$ws.Cells.Item(40, 14).Value = -12471.8  # N40 was -14271.75
$ws.Cells.Item(46, 8).Value = 1736.5  # H46 was 1858.4
$ws.Cells.Item(46, 9).Value = 2450  # I46 was 2100
$ws.Cells.Item(46, 10).Value = 1498.6666  # J46 was 1496
$ws.Cells.Item(46, 11).Value = 2450  # K46 was 2100
$ws.Cells.Item(46, 12).Value = 1498.6666  # L46 was 1496
$ws.Cells.Item(46, 13).Value = -2262  # M46 was -1912
$ws.Cells.Item(46, 14).Value = -1874.6666  # N46 was -1872
$ws.Cells.Item(126, 8).Value = 1000  # H126 was 500
$ws.Cells.Item(126, 9).Value = 1000  # I126 was 500
$ws.Cells.Item(126, 11).Value = 3000  # K126 was 1500
$ws.Cells.Item(126, 13).Value = -530  # M126 was 970
$ws.Cells.Item(132, 8).Value = 7856.7144  # H132 was 8319.6
$ws.Cells.Item(132, 9).Value = 7666.1665  # I132 was 8149.5
$ws.Cells.Item(132, 11).Value = 22998.4995  # K132 was 24448.5
$ws.Cells.Item(132, 13).Value = -20468.4995  # M132 was -21918.5
